$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Play Drago: Jewels of Fortune for Free - Expert Slot Game Review" "Play Drago: Jewels of Fortune for Free"
Replace-Text "Streak respin feature and free spins with multipliers" "Streak respin feature and free spin with multipliers and sticky symbols"
Replace-Text "High-quality graphics, with a dragon-themed design" "Crisp, high-quality graphics with rich detail"
Replace-Text "Simple and intuitive interface for desktop and mobile" "Wide range of bets for low, medium, and high stakes"
Replace-Text "Wide range of bets for low, medium, and high stakes players" "Compatible with desktop and mobile devices"
Replace-Text "High volatility can result in quickly decreasing bankroll" "Bankroll can decrease quickly during spells of volatility"
Replace-Text "Limited range of symbols compared to other slot games" "No additional bonus features apart from free spins"
Replace-Text "Read our expert review of Drago: Jewels of Fortune and play for free. Experience the dragon-themed design, streak respin feature, and free spins with multipliers." "Read our review of Drago: Jewels of Fortune and play the game for free. Discover its features and high volatility."
